# Projektstrukturplan+Budget.xlsx update
# - Projektstrukturplan: add a "NICHT AKTUELL!!!" warning label in column D, row 1
# - Projektbudget: replace the linked formulas for several phases with the
#   actual discussed ("besprochene") hours, which updates all downstream totals

$wb = $excel.ActiveWorkbook

$wsPlan = $wb.Worksheets.Item("Projektstrukturplan")
$wsBudget = $wb.Worksheets.Item("Projektbudget")

# --- Projektstrukturplan: new "NICHT AKTUELL!!!" header cell in D1 ---
$wsPlan.Range("D1").Value = "NICHT AKTUELL!!!"
$wsPlan.Range("D1").Font.Bold = $true
$wsPlan.Range("D1").Font.Size = 18
$wsPlan.Range("D1").Font.Color = 255
$wsPlan.Range("D1").Font.Name = "Calibri"

$wsPlan.Columns.Item(4).ColumnWidth = 16.5703125

# Reset the view back to the top-left corner / D1 selection
$wsPlan.Application.ActiveWindow.ScrollRow = 1
$wsPlan.Application.ActiveWindow.ScrollColumn = 1
$wsPlan.Range("D1").Select()

# --- Projektbudget: "Besprochene Zeiten" (discussed/agreed hours) entered as values ---
$wsBudget.Range("C5").Value = 115
$wsBudget.Range("C6").Value = 26
$wsBudget.Range("C7").Value = 121
$wsBudget.Range("C8").Value = 100

$wsBudget.Range("E17").Select()

$excel.Calculate()
